$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value2 = 1979
$ws.Range("I19").Value2 = 1497.75
$ws.Range("J19").Value2 = 2299.8333
$ws.Range("K19").Value2 = 1497.75
$ws.Range("L19").Value2 = 2299.8333
$ws.Range("M19").Value2 = -1322.75
$ws.Range("N19").Value2 = -2649.8333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value2 = 86283.086
$ws.Range("I64").Value2 = 127624.625
$ws.Range("J64").Value2 = 3600
$ws.Range("K64").Value2 = 127624.625
$ws.Range("L64").Value2 = 3600
$ws.Range("M64").Value2 = -127376.625
$ws.Range("N64").Value2 = -4096

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value2 = 86283.086
$ws.Range("I67").Value2 = 127624.625
$ws.Range("J67").Value2 = 3600
$ws.Range("K67").Value2 = 127624.625
$ws.Range("L67").Value2 = 3600
$ws.Range("M67").Value2 = -126766.625
$ws.Range("N67").Value2 = -5316

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value2 = 2752.5
$ws.Range("I86").Value2 = 2474.077
$ws.Range("J86").Value2 = 3081.5454
$ws.Range("K86").Value2 = 2474.077
$ws.Range("L86").Value2 = 3081.5454
$ws.Range("M86").Value2 = -1351.077
$ws.Range("N86").Value2 = -5327.5454

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value2 = 2752.5
$ws.Range("I89").Value2 = 2474.077
$ws.Range("J89").Value2 = 3081.5454
$ws.Range("K89").Value2 = 12370.385
$ws.Range("L89").Value2 = 15407.727
$ws.Range("M89").Value2 = -6754.385000000002
$ws.Range("N89").Value2 = -26639.727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value2 = 603
$ws.Range("J100").Value2 = 603
$ws.Range("L100").Value2 = 603
$ws.Range("N100").Value2 = -1685

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value2 = 2250.5
$ws.Range("J103").Value2 = 1828.5
$ws.Range("L103").Value2 = 5485.5
$ws.Range("N103").Value2 = -6657.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value2 = 3056.3333
$ws.Range("J129").Value2 = 930.25
$ws.Range("L129").Value2 = 2790.75
$ws.Range("N129").Value2 = -12790.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 68252.266
$ws.Range("J45").Value2 = 3320
$ws.Range("L45").Value2 = 3320
$ws.Range("N45").Value2 = -4074

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value2 = 20458.846
$ws.Range("J70").Value2 = 20458.846
$ws.Range("L70").Value2 = 20458.846
$ws.Range("N70").Value2 = -21044.846

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H73").Value2 = 20458.846
$ws.Range("J73").Value2 = 20458.846
$ws.Range("L73").Value2 = 20458.846
$ws.Range("N73").Value2 = -22486.846

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value2 = 1873.1538
$ws.Range("I99").Value2 = 1415
$ws.Range("J99").Value2 = 2265.8572
$ws.Range("K99").Value2 = 1415
$ws.Range("L99").Value2 = 2265.8572
$ws.Range("M99").Value2 = 83
$ws.Range("N99").Value2 = -5261.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value2 = 14612
$ws.Range("J80").Value2 = 14612
$ws.Range("L80").Value2 = 14612
$ws.Range("N80").Value2 = -16858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value2 = 14612
$ws.Range("J83").Value2 = 14612
$ws.Range("L83").Value2 = 43836
$ws.Range("N83").Value2 = -55068

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value2 = 2294.7827
$ws.Range("J99").Value2 = 2468.0625
$ws.Range("L99").Value2 = 2468.0625
$ws.Range("N99").Value2 = -5464.0625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value2 = 1081.4445
$ws.Range("I105").Value2 = 1091.4615
$ws.Range("J105").Value2 = 1055.4
$ws.Range("K105").Value2 = 1091.4615
$ws.Range("L105").Value2 = 1055.4
$ws.Range("M105").Value2 = 655.5385000000001
$ws.Range("N105").Value2 = -4549.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value2 = 2294.7827
$ws.Range("J126").Value2 = 2468.0625
$ws.Range("L126").Value2 = 7404.1875
$ws.Range("N126").Value2 = -12344.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value2 = 1600.1428
$ws.Range("I134").Value2 = 874.875
$ws.Range("K134").Value2 = 2624.625
$ws.Range("M134").Value2 = -89.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value2 = 60244
$ws.Range("J135").Value2 = 60244
$ws.Range("L135").Value2 = 60244
$ws.Range("N135").Value2 = -70384

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value2 = 4801.1665
$ws.Range("I129").Value2 = 542.2222
$ws.Range("J129").Value2 = 6626.4287
$ws.Range("K129").Value2 = 1626.6666
$ws.Range("L129").Value2 = 19879.2861
$ws.Range("M129").Value2 = 3373.3334
$ws.Range("N129").Value2 = -29879.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value2 = 796.42
$ws.Range("I131").Value2 = 468.4737
$ws.Range("J131").Value2 = 873.3457
$ws.Range("K131").Value2 = 1405.4211
$ws.Range("L131").Value2 = 2620.0371
$ws.Range("M131").Value2 = 3634.5789
$ws.Range("N131").Value2 = -12700.0371

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value2 = 24000
$ws.Range("J15").Value2 = 24000
$ws.Range("L15").Value2 = 24000
$ws.Range("N15").Value2 = -24576

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value2 = 24000
$ws.Range("J81").Value2 = 24000
$ws.Range("L81").Value2 = 24000
$ws.Range("N81").Value2 = -25996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value2 = 24000
$ws.Range("J84").Value2 = 24000
$ws.Range("L84").Value2 = 72000
$ws.Range("N84").Value2 = -81984

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value2 = 18012.75
$ws.Range("J92").Value2 = 18012.75
$ws.Range("L92").Value2 = 18012.75
$ws.Range("N92").Value2 = -21756.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value2 = 34500
$ws.Range("J93").Value2 = 34500
$ws.Range("L93").Value2 = 34500
$ws.Range("N93").Value2 = -38244

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value2 = 5000
$ws.Range("J76").Value2 = 5000
$ws.Range("L76").Value2 = 5000
$ws.Range("N76").Value2 = -5676

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H79").Value2 = 5000
$ws.Range("J79").Value2 = 5000
$ws.Range("L79").Value2 = 5000
$ws.Range("N79").Value2 = -7340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value2 = 24999.334
$ws.Range("J92").Value2 = 24999.334
$ws.Range("L92").Value2 = 24999.334
$ws.Range("N92").Value2 = -29991.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value2 = 25000
$ws.Range("J86").Value2 = 25000
$ws.Range("L86").Value2 = 25000
$ws.Range("N86").Value2 = -27246

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H89").Value2 = 25000
$ws.Range("J89").Value2 = 25000
$ws.Range("L89").Value2 = 125000
$ws.Range("N89").Value2 = -136232

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value2 = 12495
$ws.Range("J101").Value2 = 12495
$ws.Range("L101").Value2 = 12495
$ws.Range("N101").Value2 = -18985

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H102").Value2 = 43337
$ws.Range("J102").Value2 = 43337
$ws.Range("L102").Value2 = 43337
$ws.Range("N102").Value2 = -49827

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value2 = 29999.75
$ws.Range("J104").Value2 = 29999.75
$ws.Range("L104").Value2 = 29999.75
$ws.Range("N104").Value2 = -36987.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value2 = 46052
$ws.Range("J105").Value2 = 46052
$ws.Range("L105").Value2 = 46052
$ws.Range("N105").Value2 = -53040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value2 = 0
$ws.Range("J106").Value2 = 0
$ws.Range("L106").Value2 = 0
$ws.Range("N106").ClearContents()
